$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 158, shifting existing rows 158:281 down to 159:282
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row 158 with the new record
$ws.Cells.Item(158, 1).Value = 4
$ws.Cells.Item(158, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(158, 3).Value = "Los Lagos"
$ws.Cells.Item(158, 4).Value = 44741
$ws.Cells.Item(158, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(158, 5).Value = 10
$ws.Cells.Item(158, 6).Value = 100112043
$ws.Cells.Item(158, 7).Value = "Pepino ensalada"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 70
$ws.Cells.Item(158, 11).Value = 21000
$ws.Cells.Item(158, 12).Value = 21000
$ws.Cells.Item(158, 13).Value = 21000
$ws.Cells.Item(158, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(158, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(158, 16).Value = 350
$ws.Cells.Item(158, 17).Value = 60
$ws.Cells.Item(158, 18).Value = "Hortaliza"
